$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. '5.20', '8.00') need to be
# forced to Text format first, otherwise Excel auto-converts them to numbers and
# silently drops significant trailing zeros (5.20 -> 5.2, 8.00 -> 8, etc.).
$textCells = @('D5', 'D6', 'D12', 'D14', 'D19', 'D20', 'D21', 'D25', 'D26', 'D27', 'D32', 'D33', 'D34', 'D35', 'D39', 'D40', 'D41', 'D42', 'D46', 'D47', 'D48', 'D49', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.185.44'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').Value = '2.452.24'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '580.07'
$ws.Range('E5').Value = '  +1.18%  '
$ws.Range('D6').Value = '143.59'
$ws.Range('E6').Value = '  -1.79%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').Value = '2.449.65'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  -3.72%  '
$ws.Range('E11').Value = '  +2.12%  '
$ws.Range('D12').Value = '5.20'
$ws.Range('E12').Value = '  -0.98%  '
$ws.Range('E13').Value = '  -3.15%  '
$ws.Range('D14').Value = '26.53'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('E15').Value = '  -3.85%  '
$ws.Range('D16').Value = '2.832.43'
$ws.Range('E16').Value = '  -2.38%  '
$ws.Range('D17').Value = '62.135.01'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('D18').Value = '2.427.84'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').Value = '10.88'
$ws.Range('E19').Value = '  -3.69%  '
$ws.Range('D20').Value = '7.16'
$ws.Range('E20').Value = '  -2.56%  '
$ws.Range('D21').Value = '329.35'
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('E22').Value = '  -2.64%  '
$ws.Range('E23').Value = '  -3.63%  '
$ws.Range('E24').Value = '  -0.70%  '
$ws.Range('D25').Value = '65.95'
$ws.Range('E25').Value = '  +1.07%  '
$ws.Range('D26').Value = '9.37'
$ws.Range('E26').Value = '  +6.02%  '
$ws.Range('D27').Value = '619.34'
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').Value = '0.0₃0957'
$ws.Range('E28').Value = '  -6.25%  '
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('E31').Value = '  -4.75%  '
$ws.Range('D32').Value = '8.00'
$ws.Range('E32').Value = '  -2.54%  '
$ws.Range('D33').Value = '0.141'
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('D34').Value = '1.87'
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('D35').Value = '4.90'
$ws.Range('E35').Value = '  -5.51%  '
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('E37').Value = '  -6.31%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').Value = '150.30'
$ws.Range('E39').Value = '  +2.26%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').Value = '5.32'
$ws.Range('E40').Value = '  -1.53%  '
$ws.Range('D41').Value = '18.36'
$ws.Range('E41').Value = '  -2.65%  '
$ws.Range('D42').Value = '1.74'
$ws.Range('E42').Value = '  -2.67%  '
$ws.Range('E43').Value = '  +1.74%  '
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('E45').Value = '  -5.30%  '
$ws.Range('D46').Value = '143.23'
$ws.Range('E46').Value = '  -3.79%  '
$ws.Range('D47').Value = '3.64'
$ws.Range('E47').Value = '  -3.53%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.605'
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('B49').Value = 'Hedera'
$ws.Range('C49').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D49').Value = '0.0525'
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0241'
$ws.Range('E50').Value = '  +11.54%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '19.56'
$ws.Range('E51').Value = '  -7.47%  '

# Restore the original (default) cell formatting now that the text values are
# safely stored as strings, so no stray number-format styles linger on the cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
